$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "person" to "data"
$ws.Name = "data"

# New header row: idp, name, surname, personcode, user
$ws.Range("A1").Value = "idp"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "surname"
$ws.Range("D1").Value = "personcode"
$ws.Range("E1").Value = "user"

# Data rows: add idp (col A) and user (col E) values, shifting existing
# name/surname/personcode into B/C/D
$data = @(
    @{ Row = 2; Idp = "6"; Name = "John";   Surname = "Doe";       Code = "123456-89012"; User = "lv.venta.models.users.User@1722d96b" },
    @{ Row = 3; Idp = "7"; Name = "Jane";   Surname = "Smith";     Code = "123456-89012"; User = "lv.venta.models.users.User@211496df" },
    @{ Row = 4; Idp = "1"; Name = "Karina"; Surname = "Skirmante"; Code = "121212-12121"; User = "lv.venta.models.users.User@783d5f65" },
    @{ Row = 5; Idp = "2"; Name = "Karlis"; Surname = "Immers";    Code = "121212-12123"; User = "lv.venta.models.users.User@55d764b9" },
    @{ Row = 6; Idp = "3"; Name = "Janis";  Surname = "Berzins";   Code = "211221-34567"; User = "lv.venta.models.users.User@30fcecfb" },
    @{ Row = 7; Idp = "4"; Name = "Baiba";  Surname = "Kalnina";   Code = "121256-98765"; User = "lv.venta.models.users.User@158e727" },
    @{ Row = 8; Idp = "5"; Name = "Andris"; Surname = "Ribakovs";  Code = "131256-98765"; User = "lv.venta.models.users.User@6147b2ee" }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rec.Idp
    $ws.Cells.Item($r, 2).Value = $rec.Name
    $ws.Cells.Item($r, 3).Value = $rec.Surname
    $ws.Cells.Item($r, 4).Value = $rec.Code
    $ws.Cells.Item($r, 5).Value = $rec.User
}
